# "minor typos in poster" -- fix two typos on the single poster slide.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Fix 1: "Methods" box -> "D: distortion metrc (for further constraints ...)"
#     becomes "D: distortion metric (for further constraints ...)"
#     i.e. fix "metrc" -> "metric " and drop the now-redundant leading space
#     on the following run.
$shp1 = $s.Shapes.Item("Shape 117")
$tr1 = $shp1.TextFrame.TextRange

$typo = $tr1.Find("metrc")
$typo.Text = "metric "

$rest = $tr1.Find(" (for further constraints on the encoder function)")
$rest.Text = "(for further constraints on the encoder function)"

# --- Fix 2: "The All Model entails Vggnet with 3 sets ..." becomes
#     "The All Model entails VGGNet with 3 sets ..." (capitalization fix),
#     keeping the separating space but on its own run boundary.
$shp2 = $s.Shapes.Item("TextBox 17")
$tr2 = $shp2.TextFrame.TextRange

$vgg = $tr2.Find("Vggnet")
$vgg.Text = "VGGNet"

$after = $tr2.Find("with 3 sets of 2 convolutional layers with a ")
$sep = $tr2.Characters($after.Start - 1, 1)
$sep.Text = " "
